$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$row = 56

$ws.Cells.Item($row, 1).Value = "Wil je dit even bij Koen neerleggen?"
$ws.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value = "Testmail #15: Wil je dit even bij Koen neerleggen?"
$ws.Cells.Item($row, 4).Value = "Overig"
$ws.Cells.Item($row, 5).Value = "Beste heer/mevrouw,`nBedankt voor uw e-mail. Kunt u specifieker aangeven wat u precies bij Koen neergelegd wilt hebben? Zo kan ik ervoor zorgen dat uw verzoek correct wordt afgehandeld.`nMet vriendelijke groet,`n[Naam] `nE-mailassistent van [Bedrijfsnaam]"
$ws.Cells.Item($row, 6).Value = "2025-08-05 19:53:45"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Nee"
$ws.Cells.Item($row, 9).Value = "Ja"
$ws.Cells.Item($row, 10).Value = "Nee"

# Undo the automatic row-height bump that Excel applies when a
# multi-line value is entered, so the new row matches the sheet's
# default (unmodified) row formatting.
$ws.Rows.Item($row).AutoFit()

# Extend the conditional formatting ranges to cover the newly added row.
$ws.Range("D2:D55").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D56"))
$ws.Range("G2:G55").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G56"))
$ws.Range("H2:H55").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H56"))
$ws.Range("I2:I55").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I56"))
$ws.Range("J2:J55").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J56"))

# Update the Dashboard summary count for the "Overig" category.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(3, 2).Value = 10
